# Fix some issues: extend the sheet with 6 new columns (Y:AD / 25:30)
# mirroring the existing header/data layout, and adjust the frozen-pane
# selection to the newly added columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 1): alternating "boxed" header cells (Y1, AA1, AC1)
# and "plain" header cells (Z1, AB1, AD1), continuing the numbering
# sequence 35..40.
# ---------------------------------------------------------------------
$headerValues = @{ 25 = 35; 26 = 36; 27 = 37; 28 = 38; 29 = 39; 30 = 40 }
$boxedCols = @(25, 27, 29)   # Y, AA, AC
$plainCols = @(26, 28, 30)   # Z, AB, AD

foreach ($col in $boxedCols) {
    $c = $ws.Cells.Item(1, $col)
    $c.Value = $headerValues[$col]
    $c.Font.Name = "Calibri"
    $c.Font.Size = 9
    $c.Font.Color = 0
    $c.Font.Bold = $false
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.WrapText = $true
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(7).Weight = -4138
    $c.Borders.Item(7).Color = 0
    $c.Borders.Item(10).LineStyle = 1
    $c.Borders.Item(10).Weight = -4138
    $c.Borders.Item(10).Color = 0
}

foreach ($col in $plainCols) {
    $c = $ws.Cells.Item(1, $col)
    $c.Value = $headerValues[$col]
    $ws.Cells.Item(1, 2).Copy()
    $c.PasteSpecial(-4122)
    $c.Value = $headerValues[$col]
}

# ---------------------------------------------------------------------
# Data rows (2:42): all six new columns share one uniform style -
# right aligned, wrap text, light-green fill, thin-ish medium borders
# (grey on the left, black on the right), value 0.
# ---------------------------------------------------------------------
for ($row = 2; $row -le 42; $row++) {
    for ($col = 25; $col -le 30; $col++) {
        $c = $ws.Cells.Item($row, $col)
        $c.Value = 0
        $c.Font.Name = "Calibri"
        $c.Font.Size = 9
        $c.Font.Color = 0
        $c.Font.Bold = $false
        $c.HorizontalAlignment = -4152
        $c.WrapText = $true
        $c.Interior.Color = 11788485
        $c.Borders.Item(7).LineStyle = 1
        $c.Borders.Item(7).Weight = -4138
        $c.Borders.Item(7).Color = 13421772
        $c.Borders.Item(10).LineStyle = 1
        $c.Borders.Item(10).Weight = -4138
        $c.Borders.Item(10).Color = 0
    }
}

# ---------------------------------------------------------------------
# View: move the frozen-pane scroll anchor and update the active
# selection to the newly added header columns.
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("T2"))
$ws.Range("X1:AD1").Select()
